# Auto-generated edit script updating Leve profit/price cells per scheduled data refresh
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 21 (hunk 0)
$ws.Range("H21").Value = 39970
$ws.Range("I21").Value = 39970
$ws.Range("K21").Value = 39970
$ws.Range("M21").Value = -39502
# row 23 (hunk 1)
$ws.Range("H23").Value = 39970
$ws.Range("I23").Value = 39970
$ws.Range("K23").Value = 39970
$ws.Range("M23").Value = -39736
# row 51 (hunk 2)
$ws.Range("H51").Value = 9340
$ws.Range("J51").Value = 10501
$ws.Range("L51").Value = 10501
$ws.Range("N51").Value = -11469
# row 100 (hunk 3)
$ws.Range("H100").Value = 5909.8335
$ws.Range("I100").Value = 1798.2307
$ws.Range("J100").Value = 16600
$ws.Range("K100").Value = 1798.2307
$ws.Range("L100").Value = 16600
$ws.Range("M100").Value = -1257.2307
$ws.Range("N100").Value = -17682
# row 133 (hunk 4)
$ws.Range("H133").Value = 61365.547
$ws.Range("J133").Value = 61365.547
$ws.Range("L133").Value = 61365.547
$ws.Range("N133").Value = -71485.54699999999
# row 138 (hunk 5)
$ws.Range("H138").Value = 5519.467
$ws.Range("I138").Value = 3618.2856
$ws.Range("J138").Value = 6098.087
$ws.Range("K138").Value = 10854.8568
$ws.Range("L138").Value = 18294.261
$ws.Range("M138").Value = -5714.856800000001
$ws.Range("N138").Value = -28574.261

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 63 (hunk 6)
$ws.Range("H63").Value = 2896
$ws.Range("I63").Value = 2895.2
$ws.Range("J63").Value = 2900
$ws.Range("K63").Value = 2895.2
$ws.Range("L63").Value = 2900
$ws.Range("M63").Value = -2209.2
$ws.Range("N63").Value = -4272
# row 66 (hunk 7)
$ws.Range("H66").Value = 2896
$ws.Range("I66").Value = 2895.2
$ws.Range("J66").Value = 2900
$ws.Range("K66").Value = 14476
$ws.Range("L66").Value = 14500
$ws.Range("M66").Value = -11044
$ws.Range("N66").Value = -21364
# row 122 (hunk 8)
$ws.Range("H122").Value = 4269.033
$ws.Range("I122").Value = 2906.9443
$ws.Range("K122").Value = 8720.832900000001
$ws.Range("M122").Value = -6270.832900000001
# row 125 (hunk 9)
$ws.Range("H125").Value = 53763.6
$ws.Range("J125").Value = 53763.6
$ws.Range("L125").Value = 53763.6
$ws.Range("N125").Value = -63603.6
# row 132 (hunk 10)
$ws.Range("H132").Value = 6536.161
$ws.Range("I132").Value = 2035.3334
$ws.Range("K132").Value = 6106.0002
$ws.Range("M132").Value = -3576.0002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 102 (hunk 11)
$ws.Range("H102").Value = 5547
$ws.Range("I102").Value = 5547
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5547
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -2302
$ws.Range("N102").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 99 (hunk 12)
$ws.Range("H99").Value = 6201.375
$ws.Range("J99").Value = 6257
$ws.Range("L99").Value = 6257
$ws.Range("N99").Value = -9253
# row 107 (hunk 13)
$ws.Range("H107").Value = 406.36
$ws.Range("I107").Value = 266.52942
$ws.Range("K107").Value = 266.52942
$ws.Range("M107").Value = 1653.47058
# row 126 (hunk 14)
$ws.Range("H126").Value = 6201.375
$ws.Range("J126").Value = 6257
$ws.Range("L126").Value = 18771
$ws.Range("N126").Value = -23711
# row 132 (hunk 15)
$ws.Range("H132").Value = 3500.2896
$ws.Range("I132").Value = 2173.423
$ws.Range("K132").Value = 6520.268999999999
$ws.Range("M132").Value = -3990.268999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 4 (hunk 16)
$ws.Range("H4").Value = 26685736
$ws.Range("I4").Value = 1719839.6
$ws.Range("K4").Value = 5159518.800000001
$ws.Range("M4").Value = -5159406.800000001
# row 68 (hunk 17)
$ws.Range("H68").Value = 2001398.4
$ws.Range("I68").Value = 1747.25
$ws.Range("J68").Value = 10000003
$ws.Range("K68").Value = 5241.75
$ws.Range("L68").Value = 30000009
$ws.Range("M68").Value = -4430.75
$ws.Range("N68").Value = -30001631
# row 71 (hunk 18)
$ws.Range("H71").Value = 2001398.4
$ws.Range("I71").Value = 1747.25
$ws.Range("J71").Value = 10000003
$ws.Range("K71").Value = 15725.25
$ws.Range("L71").Value = 90000027
$ws.Range("M71").Value = -11669.25
$ws.Range("N71").Value = -90008139
# row 80 (hunk 19)
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# row 83 (hunk 20)
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 58 (hunk 21)
$ws.Range("H58").Value = 33046
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
# row 122 (hunk 22)
$ws.Range("H122").Value = 7783
$ws.Range("I122").Value = 6984.722
$ws.Range("J122").Value = 9579.125
$ws.Range("K122").Value = 20954.166
$ws.Range("L122").Value = 28737.375
$ws.Range("M122").Value = -18504.166
$ws.Range("N122").Value = -33637.375
# row 132 (hunk 23)
$ws.Range("H132").Value = 421178.16
$ws.Range("I132").Value = 480275.1
$ws.Range("K132").Value = 1440825.3
$ws.Range("M132").Value = -1438295.3

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7 (hunk 24)
$ws.Range("H7").Value = 1258623.2
$ws.Range("I7").Value = 1258623.2
$ws.Range("K7").Value = 1258623.2
$ws.Range("M7").Value = -1258511.2
# row 22 (hunk 25)
$ws.Range("H22").Value = 1754.6111
$ws.Range("I22").Value = 1390.2727
$ws.Range("J22").Value = 2327.1428
$ws.Range("K22").Value = 1390.2727
$ws.Range("L22").Value = 2327.1428
$ws.Range("M22").Value = -1095.2727
$ws.Range("N22").Value = -2917.1428
# row 27 (hunk 26)
$ws.Range("H27").Value = 1754.6111
$ws.Range("I27").Value = 1390.2727
$ws.Range("J27").Value = 2327.1428
$ws.Range("K27").Value = 1390.2727
$ws.Range("L27").Value = 2327.1428
$ws.Range("M27").Value = -1283.2727
$ws.Range("N27").Value = -2541.1428
# row 40 (hunk 27)
$ws.Range("H40").Value = 387357.12
$ws.Range("I40").Value = 402595.44
$ws.Range("K40").Value = 402595.44
$ws.Range("M40").Value = -402459.44
# row 64 (hunk 28)
$ws.Range("H64").Value = 30150
$ws.Range("J64").Value = 30150
$ws.Range("L64").Value = 30150
$ws.Range("N64").Value = -30600
# row 67 (hunk 29)
$ws.Range("H67").Value = 30150
$ws.Range("J67").Value = 30150
$ws.Range("L67").Value = 30150
$ws.Range("N67").Value = -31710
# row 126 (hunk 30)
$ws.Range("H126").Value = 1258623.2
$ws.Range("I126").Value = 1258623.2
$ws.Range("K126").Value = 3775869.6
$ws.Range("M126").Value = -3773399.6
# row 132 (hunk 31)
$ws.Range("H132").Value = 4543.6875
$ws.Range("I132").Value = 2671.2856
$ws.Range("K132").Value = 8013.8568
$ws.Range("M132").Value = -5483.8568

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 61 (hunk 32)
$ws.Range("H61").Value = 3906.5
$ws.Range("I61").Value = 3687.8
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 3687.8
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -3395.8
$ws.Range("N61").Value = -5584
# row 107 (hunk 33)
$ws.Range("H107").Value = 43587.543
$ws.Range("I107").Value = 54667.05
$ws.Range("J107").Value = 1485.4
$ws.Range("K107").Value = 164001.15
$ws.Range("L107").Value = 4456.200000000001
$ws.Range("M107").Value = -162081.15
$ws.Range("N107").Value = -8296.200000000001
# row 126 (hunk 34)
$ws.Range("H126").Value = 4750
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 5666.6665
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 16999.9995
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -21939.9995
